$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for B2:G11 (ME, MAE, MSE, RMSE, SE, N)
$data = @(
    @(0.02539962882261295, 0.2597273962171842, 0.3180567948316633, 0.5639652425740999, 0.5688896304229624, 52),
    @(0.4459617371143483, 0.6540342364242132, 1.388895812590984, 1.17851423945194, 1.101732372426877, 51),
    @(0.6024413171113788, 0.9796726598280953, 3.63049943662146, 1.905386951939543, 1.82599257214512, 50),
    @(0.5121426763354459, 1.006154916801229, 3.970987900588087, 1.99273377564292, 1.945754947762022, 49),
    @(0.4218352939147489, 0.8877411355809319, 3.626482144461433, 1.904332466892646, 1.876675320794537, 48),
    @(0.3949550873240204, 0.8680554327475212, 4.006566611836582, 2.001640979755506, 1.987940659805524, 39),
    @(0.3899578886789405, 0.8693687210883813, 4.090339984385427, 2.022458895598481, 2.011147094482062, 38),
    @(0.3439843485875357, 1.183619751688364, 6.937012254809386, 2.633820847136226, 2.675746881502995, 21),
    @(-0.2090860748678977, 0.8864813411966667, 3.164946506954499, 1.779029653197074, 1.833391419946016, 14),
    @(0.528408343223261, 0.528408343223261, 0.3628698470762892, 0.6023867919171944, 0.3233698924767455, 5)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = 2 + $j  # B=2 .. G=7
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}
